$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 25.57000000000056
$ws.Range("K2").Value = 58.12255376144019
$ws.Range("L2").Value = "[53.604596449494046, 62.64051107338634]"
$ws.Range("O2").Value = 1.553500271144502
$ws.Range("P2").Value = "[1.478026573760963, 1.628973968528041]"
$ws.Range("S2").Value = 54.50167644424064
$ws.Range("T2").Value = "[51.457545266856904, 57.545807621624384]"
$ws.Range("W2").Value = 19.24788788788831
$ws.Range("X2").Value = 18.94074074074116
$ws.Range("Y2").Value = 19.55503503503547

# Row 3
$ws.Range("E3").Value = 25.48000000000054
$ws.Range("K3").Value = 57.86984691193629
$ws.Range("L3").Value = "[52.82706465406821, 62.912629169804376]"
$ws.Range("O3").Value = 1.276763380738194
$ws.Range("P3").Value = "[1.1887107337907326, 1.364816027685655]"
$ws.Range("S3").Value = 53.5828118785227
$ws.Range("T3").Value = "[50.65724641341785, 56.50837734362756]"
$ws.Range("W3").Value = 20.30238238238282
$ws.Range("X3").Value = 19.94530530530573
$ws.Range("Y3").Value = 20.6594594594599
